# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp text (row 1 title) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 07:20"

# --- Swap the country labels for Groenlandia / Islas Malvinas (rows 210/211) ---
# Before: A210 = Islas Malvinas, A211 = Groenlandia
# After:  A210 = Groenlandia,    A211 = Islas Malvinas
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"

# --- Update statistics rows ---

# Row 15: Pakistan
$ws.Range("B15").Value = 270400
$ws.Range("C15").Value = 1209
$ws.Range("D15").Value = 219783
$ws.Range("E15").Value = 44854
$ws.Range("G15").Value = 54
$ws.Range("H15").Value = 5763

# Row 31: Kazajistan
$ws.Range("D31").Value = 49488
$ws.Range("E31").Value = 28413

# Row 56: Kirguistan
$ws.Range("B56").Value = 31247
$ws.Range("C56").Value = 898
$ws.Range("D56").Value = 18038
$ws.Range("E56").Value = 11998
$ws.Range("G56").Value = 42
$ws.Range("H56").Value = 1211

# Row 65: Uzbekistan
$ws.Range("B65").Value = 18986
$ws.Range("C65").Value = 118
$ws.Range("E65").Value = 8734
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 103

# Row 107: Tailandia
$ws.Range("B107").Value = 3279
$ws.Range("C107").Value = 10
$ws.Range("D107").Value = 3107
$ws.Range("E107").Value = 114
